$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 3)
$ws.Range("A3").Value = "1129.Shortest Path with Alternating Colors"
$ws.Range("B3").Value = "Medium"
$ws.Range("C3").Value = "Standard Traversal"
$ws.Range("D3").Value = "Graph BFS but with alternating colors and color tracking."
$ws.Range("E3").Value = "https://leetcode.com/problems/shortest-path-with-alternating-colors/solutions/339964/java-python-bfs/?envType=study-plan-v2&envId=graph-theory "

# Match formatting of row 2 for the new row's B and E cells
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("C3").Style = $ws.Range("C2").Style

# Add hyperlink for E3, mirroring E2's hyperlink style
$ws.Hyperlinks.Add($ws.Range("E3"), "https://leetcode.com/problems/shortest-path-with-alternating-colors/solutions/339964/java-python-bfs/?envType=study-plan-v2&envId=graph-theory ", "", "", $ws.Range("E3").Value) | Out-Null

# Adjust column widths
$ws.Columns.Item(1).ColumnWidth = 39.28515625
$ws.Columns.Item(4).ColumnWidth = 81.85546875

# Update sheet view: top-left cell and selection
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("E15").Select() | Out-Null
